# Weekly update: insert a new price record as row 199 (Feria Lagunitas de
# Puerto Montt - Ajo), pushing the existing rows 199:277 down to 200:278.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 199, shifting rows 199:277
# down to 200:278 (old row 277 becomes the new row 278).
$ws.Rows(199).Insert()

# Populate the newly inserted row 199 with the new weekly record.
$ws.Range("A199").Value = 4
$ws.Range("B199").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C199").Value = "Los Lagos"
$ws.Range("D199").Value = 44726
$ws.Range("E199").Value = 10
$ws.Range("F199").Value = 100112003
$ws.Range("G199").Value = "Ajo"
$ws.Range("H199").Value = "Chino"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 220
$ws.Range("K199").Value = 23000
$ws.Range("L199").Value = 23000
$ws.Range("M199").Value = 23000
$ws.Range("N199").Value = "`$/caja 10 kilos"
$ws.Range("O199").Value = "China"
$ws.Range("P199").Value = 2300
$ws.Range("Q199").Value = 10
$ws.Range("R199").Value = "Hortaliza"
